# ============================================================================
# Adds the "Include from SNOWMED CT" worksheet (SNOMED CT medication concepts
# for Bicillin L-A) after the existing "Metadata" sheet, and bumps the
# Metadata "Date" value to reflect the later edit time.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Metadata" sheet (the only sheet in the workbook)

# ----------------------------------------------------------------------
# 1. Bump the recorded Date value on the Metadata sheet (row 8, col B)
# ----------------------------------------------------------------------
$ws.Range("B8").Value = "2024-05-15T15:51:17+10:00"

# ----------------------------------------------------------------------
# 2. Capture the existing header / body formatting from the Metadata
#    sheet so the new sheet's look & feel matches it exactly.
# ----------------------------------------------------------------------
$hdrSrc = $ws.Range("A1")
$bodySrc = $ws.Range("A2")

$hdrBold = $hdrSrc.Font.Bold
$hdrFontName = $hdrSrc.Font.Name
$hdrFontSize = $hdrSrc.Font.Size
$hdrFillColor = $hdrSrc.Interior.Color

$bodyFontName = $bodySrc.Font.Name
$bodyFontSize = $bodySrc.Font.Size

# Thin, light-grey grid border (matches the border used throughout the
# Metadata sheet) expressed as an explicit RGB so it renders identically
# regardless of which indexed-colour palette the host uses.
$borderColor = 12632256  # RGB(192,192,192)

# ----------------------------------------------------------------------
# 3. Create the new worksheet right after "Metadata"
# ----------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "Include from SNOWMED CT"

# ----------------------------------------------------------------------
# 4. Populate the cell values.
#    SNOMED CT concept ids are long numeric-looking codes - prefix them
#    with an apostrophe so Excel stores them as text instead of coercing
#    them to numbers (which would lose precision / use sci. notation).
# ----------------------------------------------------------------------
$newSheet.Range("A1").Value = "Concept"
$newSheet.Range("B1").Value = "Description"

$newSheet.Range("A2").Value = "'1276211000168102"
$newSheet.Range("B2").Value = "Bicillin L-A 1.2 million units (1016.6 mg)/2.3 mL injection, 2.3 mL syringe"

$newSheet.Range("A3").Value = "'70271000036108"
$newSheet.Range("B3").Value = "Bicillin L-A 1.2 million units (900 mg)/2 mL injection, 2 mL cartridge"

$newSheet.Range("A4").Value = "'39212011000036104"
$newSheet.Range("B4").Value = "Bicillin L-A 1.2 million units (900 mg)/2 mL injection, 2 mL syringe"

$newSheet.Range("A5").Value = "'1248021000168101"
$newSheet.Range("B5").Value = "Bicillin L-A 600 000 units (517 mg)/1.17 mL injection, 1.17 mL syringe"

$newSheet.Range("A6").Value = "'4089011000036101"
$newSheet.Range("B6").Value = "Bicillin L-A"

# Row 7 is intentionally blank (still carries the body style, no value)
$newSheet.Range("A7").Style = $newSheet.Range("A7").Style
$newSheet.Range("B7").Style = $newSheet.Range("B7").Style

$newSheet.Range("A8").Value = "System URI"
$newSheet.Range("B8").Value = "http://snomed.info/sct"

# ----------------------------------------------------------------------
# 5. Formatting - header row (bold, filled, top-aligned, wrapped, bordered)
# ----------------------------------------------------------------------
$hdrRange = $newSheet.Range("A1:B1")
$hdrRange.Font.Bold = $hdrBold
$hdrRange.Font.Name = $hdrFontName
$hdrRange.Font.Size = $hdrFontSize
$hdrRange.Interior.Color = $hdrFillColor
$hdrRange.VerticalAlignment = -4160  # xlTop
$hdrRange.WrapText = $true
$hdrRange.Borders.LineStyle = 1
$hdrRange.Borders.Color = $borderColor

# ----------------------------------------------------------------------
# 6. Formatting - body rows (top-aligned, wrapped, bordered)
# ----------------------------------------------------------------------
$bodyRange = $newSheet.Range("A2:B8")
$bodyRange.Font.Name = $bodyFontName
$bodyRange.Font.Size = $bodyFontSize
$bodyRange.VerticalAlignment = -4160  # xlTop
$bodyRange.WrapText = $true
$bodyRange.Borders.LineStyle = 1
$bodyRange.Borders.Color = $borderColor

# ----------------------------------------------------------------------
# 7. Column widths (characters) - mirrors the authored widths of
#    30.703125 / 50.703125
# ----------------------------------------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 30.7
$newSheet.Columns.Item(2).ColumnWidth = 50.7

# ----------------------------------------------------------------------
# 8. Restore "Metadata" as the active/selected sheet
# ----------------------------------------------------------------------
$ws.Activate()
